# Updated BGR model - 2025-08-13 11:26
# Swap the ordering of adjacent "cost class" rows (and their lcoe_class rank
# in column P) on the "wind" sheet for CF classes won-BGR_29, won-BGR_25,
# won-BGR_21 and won-BGR_17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("wind")

$rowPairs = @(
    @(4, 5),
    @(13, 14),
    @(27, 28),
    @(47, 48)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # Column C (process name), D (description) and K (process name, mirrored)
    foreach ($col in 3, 4, 11) {
        $v1 = $ws.Cells.Item($r1, $col).Value()
        $v2 = $ws.Cells.Item($r2, $col).Value()
        $ws.Cells.Item($r1, $col).Value = $v2
        $ws.Cells.Item($r2, $col).Value = $v1
    }

    # Column P (lcoe_class rank)
    $p1 = $ws.Cells.Item($r1, 16).Value()
    $p2 = $ws.Cells.Item($r2, 16).Value()
    $ws.Cells.Item($r1, 16).Value = $p2
    $ws.Cells.Item($r2, 16).Value = $p1
}
